# Add group to header file PCHS app
#
# The header table in Sheet1 lists "Report type" / "Notes" / (blank) labels
# in column A with their values in column B. This adds a new "Group" /
# "Site" row, and normalises the old "Microsoft Sans Serif" label font to
# Calibri (matching the rest of the workbook's theme font) along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-A labels ("Report type", "Notes") were still set in the legacy
# "Microsoft Sans Serif" bold font; bring them onto Calibri like everything
# else in the sheet.
$ws.Range("A1:A2").Font.Name = "Calibri"

# A3 (blank) used the legacy plain "Microsoft Sans Serif" font too.
$ws.Range("A3").Font.Name = "Calibri"

# New row 4: "Group" label (bold, accent-coloured like the sheet's other
# emphasised text) and its "Site" value.
$ws.Range("A4").Value = "Group"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.ThemeColor = 1

$ws.Range("B4").Value = "Site"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

# Match the author's final cursor position recorded in the workbook.
$ws.Range("D10").Select() | Out-Null
